# AlphaFiberF-HW50.xlsx -- add averaged-intensity results for the new
# "spiral" sampling schemes (Gaussian-Quadrature + 3 Spiral-* variants),
# inserted ahead of the NoRotation/Rotation/HexGrid rows, and append the
# 3 HexGrid rows that were pushed out to the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table by 3 rows: duplicate the last existing data row (16),
# formatting included, down into rows 17-19 so the new rows pick up the
# same bold/centered/bordered "index" style used by column A without
# introducing any new style definitions.
$ws.Range("A16:M16").Copy($ws.Range("A17:M17"))
$ws.Range("A16:M16").Copy($ws.Range("A18:M18"))
$ws.Range("A16:M16").Copy($ws.Range("A19:M19"))

# Fix up the running index in column A for the 3 newly-added rows.
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

# Re-label rows 10-19 (column B) with the final scheme order: the new
# Gaussian-Quadrature + Spiral-* schemes are inserted right after the
# "Ring Perpendicular to TD" row, which pushes NoRotation/Rotation/HexGrid
# further down (HexGrid now lands in the 3 brand-new rows).
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# Columns C:M for rows 17-19 already came through as 1 (copied from row
# 16), matching every other data row's "all intensities present" flags.
